# Generate Report for Handback
#
# The "Ready for handoff" status (row 3 / the 97c718ce-...md file) is now
# known to have failed its handback transform, and both locale sheets get a
# new "Error Detail" note (column K) explaining the file-name mismatch that
# caused the failure.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status for the 97c718ce-...md row flips from "Ready for handoff" to
# "Handback transform failed" everywhere it is shown.
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhcn.Range("C3").Value     = "Handback transform failed"
$dede.Range("C3").Value     = "Handback transform failed"

# New Error Detail (column K) entries describing the handback/handoff file
# name mismatch for each locale.
$zhcn.Range("K3").Value = "Handback file name: pd4t25lm.2uf is different with handoff file name: 97c718ce-fbc4-454f-b4c7-ef0700834efb.9d9fe7d578988213b37d7c7b52a7b027026123f4.zh-cn."
$dede.Range("K3").Value = "Handback file name: pd4t25lm.2uf is different with handoff file name: 97c718ce-fbc4-454f-b4c7-ef0700834efb.9d9fe7d578988213b37d7c7b52a7b027026123f4.de-de."
